# Registers.xlsx -- "added TMF8801 as an extra I2C peripheral to test"
#
# The TMF8801 register map (worksheet "TMF8801") gets:
#   - APPREV_MAJOR (row 3) narrowed from a 16-bit field to an 8-bit field
#   - a new APPREV_MINOR row inserted right after APPREQID (row 5), also 8-bit
#   - the pre-existing APPREV_MINOR / APPREV_PATCH rows (now rows 23/24)
#     narrowed from 16-bit to 8-bit fields as well
#   - ten new rows appended describing the ENABLE register's individual bits
#     (CPU_RESET/CPU_READY/PON) plus several new registers (RESULT_INFO,
#     RELIABILITY, MEAS_STATUS, SYS_CLOCK, DISTANCE_PEAK, TEMPERATURE,
#     STATE_FIELD)

$wb = $excel.ActiveWorkbook
$tmf = $wb.Worksheets.Item("TMF8801")

# --- APPREV_MAJOR (row 3) narrows from 16/15 to 8/7 --------------------
$tmf.Range("D3").Value = 8
$tmf.Range("E3").Value = 7

# --- insert a new APPREV_MINOR row right after APPREQID (row 4) --------
$tmf.Rows.Item(5).Insert()
$tmf.Range("A5").Value = "APPREV_MINOR"
$tmf.Range("B5").Value = "0x12"
$tmf.Range("C5").Value = "0x00"
$tmf.Range("D5").Value = 8
$tmf.Range("E5").Value = 7
$tmf.Range("F5").Value = 0

# --- the original APPREV_MINOR / APPREV_PATCH rows (shifted down by the
#     insert above to rows 23/24) also narrow from 16/15 to 8/7 ----------
$tmf.Range("D23").Value = 8
$tmf.Range("E23").Value = 7
$tmf.Range("D24").Value = 8
$tmf.Range("E24").Value = 7

# --- append the new bit-field / register rows at the bottom ------------
$tmf.Range("A34").Value = "CPU_RESET"
$tmf.Range("B34").Value = "0xE0"
$tmf.Range("C34").Value = "0x0"
$tmf.Range("D34").Value = 1
$tmf.Range("E34").Value = 7
$tmf.Range("F34").Value = 7

$tmf.Range("A35").Value = "CPU_READY"
$tmf.Range("B35").Value = "0xE0"
$tmf.Range("C35").Value = "0x0"
$tmf.Range("D35").Value = 1
$tmf.Range("E35").Value = 6
$tmf.Range("F35").Value = 6

$tmf.Range("A36").Value = "PON"
$tmf.Range("B36").Value = "0xE0"
$tmf.Range("C36").Value = "0x0"
$tmf.Range("D36").Value = 1
$tmf.Range("E36").Value = 0
$tmf.Range("F36").Value = 0

$tmf.Range("A37").Value = "RESULT_INFO"
$tmf.Range("B37").Value = "0x21"
$tmf.Range("C37").Value = "0x00"
$tmf.Range("D37").Value = 8
$tmf.Range("E37").Value = 7
$tmf.Range("F37").Value = 0

$tmf.Range("A38").Value = "RELIABILITY"
$tmf.Range("B38").Value = "0x21"
$tmf.Range("C38").Value = "0x00"
$tmf.Range("D38").Value = 6
$tmf.Range("E38").Value = 5
$tmf.Range("F38").Value = 0

$tmf.Range("A39").Value = "MEAS_STATUS"
$tmf.Range("B39").Value = "0x21"
$tmf.Range("C39").Value = "0x00"
$tmf.Range("D39").Value = 2
$tmf.Range("E39").Value = 7
$tmf.Range("F39").Value = 6

# SYS_CLOCK typed first (row 40) ...
$tmf.Range("A40").Value = "SYS_CLOCK"
$tmf.Range("B40").Value = "0x24"
$tmf.Range("C40").Value = "0x00"
$tmf.Range("D40").Value = 32
$tmf.Range("E40").Value = 31
$tmf.Range("F40").Value = 0

# ... then DISTANCE_PEAK gets inserted above it, pushing SYS_CLOCK to row 41
$tmf.Rows.Item(40).Insert()
$tmf.Range("A40").Value = "DISTANCE_PEAK"
$tmf.Range("B40").Value = "0x23"
$tmf.Range("C40").Value = "0x00"
$tmf.Range("D40").Value = 16
$tmf.Range("E40").Value = 15
$tmf.Range("F40").Value = 0

$tmf.Range("A42").Value = "TEMPERATURE"
$tmf.Range("B42").Value = "0x32"
$tmf.Range("C42").Value = "0x00"
$tmf.Range("D42").Value = 8
$tmf.Range("E42").Value = 7
$tmf.Range("F42").Value = 0

$tmf.Range("A43").Value = "STATE_FIELD"
$tmf.Range("B43").Value = "0x1C"
$tmf.Range("C43").Value = "0x00"
$tmf.Range("D43").Value = 8
$tmf.Range("E43").Value = 8
$tmf.Range("F43").Value = 0

# --- sheet-view bookkeeping ---------------------------------------------
# TCA9555: selection moved from A1:F5 to E6
$tca = $wb.Worksheets.Item("TCA9555")
$tca.Range("E6").Select()

# TMF8801 stays the active tab; selection ends on the newly-entered
# RESULT_INFO row (the whole row, as if selected via the row header)
$tmf.Select()
$tmf.Rows.Item(37).Select()
